# Updates crypto price values in column D (Price) to reflect the latest
# scrape, matching the commit "Updated symbol list on Mon Dec 19 05:22:24
# UTC 2022 with GitHub Actions".
#
# Column D cells are stored as text (inline strings) in the source
# workbook, e.g. "247.70", preserving trailing zeros / significant
# digits as scraped. Prefixing the assigned value with a single quote
# forces Excel to keep it as text instead of silently coercing it to a
# numeric value (which would lose formatting like trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Value = "247.70" },
    @{ Row = 3; Value = "21.80" },
    @{ Row = 4; Value = "5.515" },
    @{ Row = 5; Value = "0.05651" },
    @{ Row = 6; Value = "3.381" },
    @{ Row = 7; Value = "6.437" },
    @{ Row = 8; Value = "0.8023" },
    @{ Row = 9; Value = "1.034" },
    @{ Row = 10; Value = "0.1423" },
    @{ Row = 11; Value = "0.07243" },
    @{ Row = 12; Value = "0.03140" },
    @{ Row = 13; Value = "0.02949" },
    @{ Row = 14; Value = "0.09284" },
    @{ Row = 15; Value = "0.001654" },
    @{ Row = 16; Value = "3.216" },
    @{ Row = 17; Value = "0.04745" },
    @{ Row = 18; Value = "0.0005857" },
    @{ Row = 19; Value = "0.006408" },
    @{ Row = 20; Value = "0.005020" },
    @{ Row = 22; Value = "0.0001502" },
    @{ Row = 23; Value = "0.0003204" },
    @{ Row = 24; Value = "4.033" },
    @{ Row = 25; Value = "2.110" },
    @{ Row = 26; Value = "0.3270" },
    @{ Row = 40; Value = "0.04091" },
    @{ Row = 41; Value = "0.006919" },
    @{ Row = 42; Value = "0.1043" },
    @{ Row = 43; Value = "0.002974" },
    @{ Row = 44; Value = "0.009120" },
    @{ Row = 45; Value = "0.00005827" },
    @{ Row = 47; Value = "0.7861" },
    @{ Row = 48; Value = "0.01595" },
    @{ Row = 49; Value = "0.00002103" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = "'" + $u.Value
}
